$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 77 (shifts existing rows 77-200 down to 78-201)
$ws.Rows.Item(77).Insert()

# Populate the newly inserted row 77 with the new data record
$ws.Cells.Item(77, 1).Value = 3
$ws.Cells.Item(77, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(77, 3).Value = "Coquimbo"
$ws.Cells.Item(77, 4).Value = 44477
$ws.Cells.Item(77, 5).Value = 5
$ws.Cells.Item(77, 6).Value = 100114013
$ws.Cells.Item(77, 7).Value = "Zanahoria"
$ws.Cells.Item(77, 8).Value = "Sin especificar"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value = 340
$ws.Cells.Item(77, 11).Value = 8000
$ws.Cells.Item(77, 12).Value = 8500
$ws.Cells.Item(77, 13).Value = 8235
$ws.Cells.Item(77, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(77, 15).Value = "Chillán"
$ws.Cells.Item(77, 16).Value = 412
$ws.Cells.Item(77, 17).Value = 20
$ws.Cells.Item(77, 18).Value = "Hortaliza"
